# Auto-generated: applies scraped market-data refresh to leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 69506.39999999999
$ws.Cells.Item(64, 9).Value = 145442.58
$ws.Cells.Item(64, 11).Value = 145442.58
$ws.Cells.Item(64, 13).Value = -145194.58
$ws.Cells.Item(67, 8).Value = 69506.39999999999
$ws.Cells.Item(67, 9).Value = 145442.58
$ws.Cells.Item(67, 11).Value = 145442.58
$ws.Cells.Item(67, 13).Value = -144584.58
$ws.Cells.Item(76, 8).Value = 2766.8667
$ws.Cells.Item(76, 9).Value = 2766.8667
$ws.Cells.Item(76, 11).Value = 2766.8667
$ws.Cells.Item(76, 13).Value = -2451.8667
$ws.Cells.Item(79, 8).Value = 2766.8667
$ws.Cells.Item(79, 9).Value = 2766.8667
$ws.Cells.Item(79, 11).Value = 2766.8667
$ws.Cells.Item(79, 13).Value = -1674.8667
$ws.Cells.Item(87, 8).Value = 0
$ws.Cells.Item(87, 10).Value = 0
$ws.Cells.Item(87, 12).ClearContents()
$ws.Cells.Item(87, 14).Value = 0
$ws.Cells.Item(90, 8).Value = 0
$ws.Cells.Item(90, 10).Value = 0
$ws.Cells.Item(90, 12).ClearContents()
$ws.Cells.Item(90, 14).Value = 0
$ws.Cells.Item(132, 8).Value = 22023.783
$ws.Cells.Item(132, 9).Value = 3302.6592
$ws.Cells.Item(132, 11).Value = 9907.9776
$ws.Cells.Item(132, 13).Value = -7377.9776
$ws.Cells.Item(137, 8).Value = 2775.5
$ws.Cells.Item(137, 9).Value = 1104.1818
$ws.Cells.Item(137, 10).Value = 3316.2207
$ws.Cells.Item(137, 11).Value = 3312.5454
$ws.Cells.Item(137, 12).Value = 9948.6621
$ws.Cells.Item(137, 13).Value = -762.5454
$ws.Cells.Item(137, 14).Value = -15048.6621
$ws.Cells.Item(138, 8).Value = 2881.1912
$ws.Cells.Item(138, 9).Value = 2577
$ws.Cells.Item(138, 11).Value = 7731
$ws.Cells.Item(138, 13).Value = -2591

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 27508
$ws.Cells.Item(32, 9).Value = 27985.852
$ws.Cells.Item(32, 10).Value = 11500
$ws.Cells.Item(32, 11).Value = 27985.852
$ws.Cells.Item(32, 12).Value = 11500
$ws.Cells.Item(32, 13).Value = -27698.852
$ws.Cells.Item(32, 14).Value = -12074
$ws.Cells.Item(63, 8).Value = 2761.25
$ws.Cells.Item(63, 9).Value = 2395.5264
$ws.Cells.Item(63, 10).Value = 3533.3333
$ws.Cells.Item(63, 11).Value = 2395.5264
$ws.Cells.Item(63, 12).Value = 3533.3333
$ws.Cells.Item(63, 13).Value = -1709.5264
$ws.Cells.Item(63, 14).Value = -4905.3333
$ws.Cells.Item(66, 8).Value = 2761.25
$ws.Cells.Item(66, 9).Value = 2395.5264
$ws.Cells.Item(66, 10).Value = 3533.3333
$ws.Cells.Item(66, 11).Value = 11977.632
$ws.Cells.Item(66, 12).Value = 17666.6665
$ws.Cells.Item(66, 13).Value = -8545.632000000001
$ws.Cells.Item(66, 14).Value = -24530.6665
$ws.Cells.Item(80, 8).Value = 61650
$ws.Cells.Item(80, 10).Value = 61650
$ws.Cells.Item(80, 12).Value = 61650
$ws.Cells.Item(80, 14).Value = -63646
$ws.Cells.Item(83, 8).Value = 61650
$ws.Cells.Item(83, 10).Value = 61650
$ws.Cells.Item(83, 12).Value = 184950
$ws.Cells.Item(83, 14).Value = -194934

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2241.0417
$ws.Cells.Item(105, 9).Value = 2127.4
$ws.Cells.Item(105, 10).Value = 2322.2144
$ws.Cells.Item(105, 11).Value = 2127.4
$ws.Cells.Item(105, 12).Value = 2322.2144
$ws.Cells.Item(105, 13).Value = -380.4000000000001
$ws.Cells.Item(105, 14).Value = -5816.2144

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 2696.6191
$ws.Cells.Item(62, 9).Value = 2575.2666
$ws.Cells.Item(62, 10).Value = 3000
$ws.Cells.Item(62, 11).Value = 2575.2666
$ws.Cells.Item(62, 12).Value = 3000
$ws.Cells.Item(62, 13).Value = -1951.2666
$ws.Cells.Item(62, 14).Value = -4248
$ws.Cells.Item(65, 8).Value = 2696.6191
$ws.Cells.Item(65, 9).Value = 2575.2666
$ws.Cells.Item(65, 10).Value = 3000
$ws.Cells.Item(65, 11).Value = 12876.333
$ws.Cells.Item(65, 12).Value = 15000
$ws.Cells.Item(65, 13).Value = -9756.332999999999
$ws.Cells.Item(65, 14).Value = -21240
$ws.Cells.Item(81, 8).Value = 39996
$ws.Cells.Item(81, 10).Value = 39996
$ws.Cells.Item(81, 12).Value = 39996
$ws.Cells.Item(81, 14).Value = -41992
$ws.Cells.Item(84, 8).Value = 39996
$ws.Cells.Item(84, 10).Value = 39996
$ws.Cells.Item(84, 12).Value = 119988
$ws.Cells.Item(84, 14).Value = -129972
$ws.Cells.Item(88, 8).Value = 43735.4
$ws.Cells.Item(88, 10).Value = 43735.4
$ws.Cells.Item(88, 12).Value = 43735.4
$ws.Cells.Item(88, 14).Value = -44547.4
$ws.Cells.Item(91, 8).Value = 43735.4
$ws.Cells.Item(91, 10).Value = 43735.4
$ws.Cells.Item(91, 12).Value = 43735.4
$ws.Cells.Item(91, 14).Value = -46543.4
$ws.Cells.Item(125, 8).Value = 49326
$ws.Cells.Item(125, 10).Value = 49326
$ws.Cells.Item(125, 12).Value = 49326
$ws.Cells.Item(125, 14).Value = -54246

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4591.206
$ws.Cells.Item(70, 9).Value = 4447.421
$ws.Cells.Item(70, 10).Value = 4773.3335
$ws.Cells.Item(70, 11).Value = 4447.421
$ws.Cells.Item(70, 12).Value = 4773.3335
$ws.Cells.Item(70, 13).Value = -4177.421
$ws.Cells.Item(70, 14).Value = -5313.3335
$ws.Cells.Item(73, 8).Value = 4591.206
$ws.Cells.Item(73, 9).Value = 4447.421
$ws.Cells.Item(73, 10).Value = 4773.3335
$ws.Cells.Item(73, 11).Value = 4447.421
$ws.Cells.Item(73, 12).Value = 4773.3335
$ws.Cells.Item(73, 13).Value = -3511.421
$ws.Cells.Item(73, 14).Value = -6645.3335
$ws.Cells.Item(80, 8).Value = 195884.84
$ws.Cells.Item(80, 9).Value = 337286.66
$ws.Cells.Item(80, 10).Value = 3064.182
$ws.Cells.Item(80, 11).Value = 337286.66
$ws.Cells.Item(80, 12).Value = 3064.182
$ws.Cells.Item(80, 13).Value = -336288.66
$ws.Cells.Item(80, 14).Value = -5060.182
$ws.Cells.Item(83, 8).Value = 195884.84
$ws.Cells.Item(83, 9).Value = 337286.66
$ws.Cells.Item(83, 10).Value = 3064.182
$ws.Cells.Item(83, 11).Value = 1686433.3
$ws.Cells.Item(83, 12).Value = 15320.91
$ws.Cells.Item(83, 13).Value = -1681441.3
$ws.Cells.Item(83, 14).Value = -25304.91
$ws.Cells.Item(118, 8).Value = 35573.5
$ws.Cells.Item(118, 10).Value = 35573.5
$ws.Cells.Item(118, 12).Value = 35573.5
$ws.Cells.Item(118, 14).Value = -38887.5
$ws.Cells.Item(120, 8).Value = 25658.5
$ws.Cells.Item(120, 10).Value = 25658.5
$ws.Cells.Item(120, 12).Value = 25658.5
$ws.Cells.Item(120, 14).Value = -35334.5
$ws.Cells.Item(127, 8).Value = 42217.332
$ws.Cells.Item(127, 10).Value = 42217.332
$ws.Cells.Item(127, 12).Value = 42217.332
$ws.Cells.Item(127, 14).Value = -52137.332
$ws.Cells.Item(131, 8).Value = 37995
$ws.Cells.Item(131, 10).Value = 37995
$ws.Cells.Item(131, 12).Value = 37995
$ws.Cells.Item(131, 14).Value = -48075

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(88, 8).Value = 24608.25
$ws.Cells.Item(88, 10).Value = 44181
$ws.Cells.Item(88, 12).Value = 44181
$ws.Cells.Item(88, 14).Value = -45037
$ws.Cells.Item(91, 8).Value = 24608.25
$ws.Cells.Item(91, 10).Value = 44181
$ws.Cells.Item(91, 12).Value = 44181
$ws.Cells.Item(91, 14).Value = -47145
$ws.Cells.Item(92, 8).Value = 40389
$ws.Cells.Item(92, 10).Value = 40389
$ws.Cells.Item(92, 12).Value = 40389
$ws.Cells.Item(92, 14).Value = -45381
$ws.Cells.Item(96, 8).Value = 31547.25
$ws.Cells.Item(96, 10).Value = 31547.25
$ws.Cells.Item(96, 12).Value = 31547.25
$ws.Cells.Item(96, 14).Value = -37039.25
$ws.Cells.Item(99, 8).Value = 25876.5
$ws.Cells.Item(99, 9).Value = 17629.5
$ws.Cells.Item(99, 11).Value = 17629.5
$ws.Cells.Item(99, 13).Value = -14634.5
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 12).ClearContents()
$ws.Cells.Item(102, 14).Value = 0
$ws.Cells.Item(117, 8).Value = 31326
$ws.Cells.Item(117, 10).Value = 31326
$ws.Cells.Item(117, 12).Value = 31326
$ws.Cells.Item(117, 14).Value = -40504
$ws.Cells.Item(123, 8).Value = 35631.332
$ws.Cells.Item(123, 10).Value = 35631.332
$ws.Cells.Item(123, 12).Value = 35631.332
$ws.Cells.Item(123, 14).Value = -45431.332
$ws.Cells.Item(129, 8).Value = 38743.332
$ws.Cells.Item(129, 10).Value = 38743.332
$ws.Cells.Item(129, 12).Value = 38743.332
$ws.Cells.Item(129, 14).Value = -48743.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(27, 8).Value = 31536.334
$ws.Cells.Item(27, 10).Value = 31536.334
$ws.Cells.Item(27, 12).Value = 31536.334
$ws.Cells.Item(27, 14).Value = -31674.334
$ws.Cells.Item(93, 8).Value = 35571.43
$ws.Cells.Item(93, 10).Value = 35571.43
$ws.Cells.Item(93, 12).Value = 35571.43
$ws.Cells.Item(93, 14).Value = -40563.43
$ws.Cells.Item(97, 8).Value = 39072
$ws.Cells.Item(97, 10).Value = 39072
$ws.Cells.Item(97, 12).Value = 39072
$ws.Cells.Item(97, 14).Value = -41054
$ws.Cells.Item(109, 8).Value = 37369
$ws.Cells.Item(109, 10).Value = 37369
$ws.Cells.Item(109, 12).Value = 37369
$ws.Cells.Item(109, 14).Value = -40143
$ws.Cells.Item(115, 8).Value = 30303.75
$ws.Cells.Item(115, 10).Value = 30303.75
$ws.Cells.Item(115, 12).Value = 30303.75
$ws.Cells.Item(115, 14).Value = -33437.75
$ws.Cells.Item(118, 8).Value = 48963.6
$ws.Cells.Item(118, 10).Value = 48963.6
$ws.Cells.Item(118, 12).Value = 48963.6
$ws.Cells.Item(118, 14).Value = -52277.6
$ws.Cells.Item(129, 8).Value = 43429
$ws.Cells.Item(129, 9).Value = 0
$ws.Cells.Item(129, 11).Value = 0
$ws.Cells.Item(129, 13).ClearContents()
$ws.Cells.Item(140, 8).Value = 49329.42
$ws.Cells.Item(140, 10).Value = 50181.055
$ws.Cells.Item(140, 12).Value = 50181.055
$ws.Cells.Item(140, 14).Value = -60541.055
